# Add two new shared/payment-provider rows (37-39) to the netCrypto ledger
# on SheetName1, matching new deposit/withdrawal transactions, and update
# the view/selection state to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New transaction rows -------------------------------------------------

# Row 37: Deposit / Wiretransfer / Roobic
$ws.Range("E37").Value = "Deposit"
$ws.Range("N37").Value = "Wiretransfer"
$ws.Range("P37").Value = "Roobic"
$ws.Range("T37").Value = 2541.0567000000001

# Row 38: Withdrawal / Credit Card / Stratrading
$ws.Range("E38").Value = "Withdrawal"
$ws.Range("N38").Value = "Credit Card"
$ws.Range("P38").Value = "Stratrading"
$ws.Range("T38").Value = 271.51499999999999

# Row 39: Withdrawal / Crypto / USDT ERC20
$ws.Range("E39").Value = "Withdrawal"
$ws.Range("N39").Value = "Crypto"
$ws.Range("P39").Value = "USDT ERC20"
$ws.Range("T39").Value = 271.91500000000002

# --- View / selection state -------------------------------------------------

$ws.Activate()

# Scroll the sheet so row 15 is at the top of the viewport.
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1

# Leave the final selection on K46, matching where the user clicked next.
[void]$ws.Range("K46").Select()
